$wb = $excel.ActiveWorkbook

$newTime = "02:57:00"

# ---------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: $newTime"
$ws1.Range("A3").Value = "Total filas: 5"

# Row 6 (existing) - update Hora_Scrap and Minutos
$ws1.Range("A6").Value = $newTime
$ws1.Range("B6").Value = "02:58"
$ws1.Range("C6").Value = "215_ALUAR"
$ws1.Range("D6").Value = 1
$ws1.Range("E6").Value = "LP1912"

# Row 7 (existing) - update Hora_Scrap, Hora_Llegada, Minutos
$ws1.Range("A7").Value = $newTime
$ws1.Range("B7").Value = "03:48"
$ws1.Range("C7").Value = "14_ABASTO"
$ws1.Range("D7").Value = 51
$ws1.Range("E7").Value = "LP1912"

# Row 8 (existing) - update Hora_Scrap and Minutos
$ws1.Range("A8").Value = $newTime
$ws1.Range("B8").Value = "04:01"
$ws1.Range("C8").Value = "81_EL PELIGRO"
$ws1.Range("D8").Value = 64
$ws1.Range("E8").Value = "LP1912"

# Row 9 (new)
$ws1.Range("A9").Value = $newTime
$ws1.Range("B9").Value = "04:46"
$ws1.Range("C9").Value = "215A_EL PATO"
$ws1.Range("D9").Value = 109
$ws1.Range("E9").Value = "LP1912"

# Row 10 (new)
$ws1.Range("A10").Value = $newTime
$ws1.Range("B10").Value = "04:53"
$ws1.Range("C10").Value = "11_ETCHEVERRY"
$ws1.Range("D10").Value = 116
$ws1.Range("E10").Value = "LP1912"

# ---------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: $newTime"
$ws2.Range("A3").Value = "Total filas: 2"

# Row 6 (existing) - update Hora_Scrap and Minutos
$ws2.Range("A6").Value = $newTime
$ws2.Range("B6").Value = "02:58"
$ws2.Range("C6").Value = "215_ALUAR"
$ws2.Range("D6").Value = 1
$ws2.Range("E6").Value = "LP1912"

# Row 7 (new)
$ws2.Range("A7").Value = $newTime
$ws2.Range("B7").Value = "04:46"
$ws2.Range("C7").Value = "215A_EL PATO"
$ws2.Range("D7").Value = 109
$ws2.Range("E7").Value = "LP1912"

# ---------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: $newTime"
